$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ELC_AVA")

# Insert a new column at J (shifts the old Pset_PN column from J to K)
$ws.Columns.Item(10).Insert()

# The freshly inserted column picks up leftover styled-but-empty cells in
# rows that only had a style (no value) in the old J column (rows 9-11).
# Clean those back up so the column stays untouched there, same as the
# source rows which never had data past column I.
$ws.Range("J9:J11").ClearFormats()
$ws.Range("J9:J11").ClearContents()

# New header for the inserted "MAR" region column
$ws.Range("J4").Value = "MAR"

# PV is not available for the existing techs/region combos, only for the
# brand new MAR-only solar PV tech added below
$ws.Range("J5").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("J7").Value = 0

# New row for the solar PV technology, available only in MAR
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = "ERSOLPV5N"
